$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = " Weight MeOH_2"
$ws.Range("E3").Value = " Weight Zn6H_2"
